$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 119; this shifts the existing rows 119-126
# down to 120-127 (formats/styles carried along automatically).
$ws.Rows.Item(119).Insert()

# Fill in the brand-new row 119 with its data.
$ws.Cells.Item(119, 1).Value = 4
$ws.Cells.Item(119, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(119, 3).Value = "Los Lagos"
$ws.Cells.Item(119, 4).Value = 44491
$ws.Cells.Item(119, 5).Value = 10
$ws.Cells.Item(119, 6).Value = 100112028
$ws.Cells.Item(119, 7).Value = "Sandia"
$ws.Cells.Item(119, 8).Value = "Sin especificar"
$ws.Cells.Item(119, 9).Value = "Primera"
$ws.Cells.Item(119, 10).Value = 600
$ws.Cells.Item(119, 11).Value = 1000
$ws.Cells.Item(119, 12).Value = 1000
$ws.Cells.Item(119, 13).Value = 1000
$ws.Cells.Item(119, 14).Value = "`$/kilo (volumen en unidades)"
$ws.Cells.Item(119, 15).Value = "Perú"
$ws.Cells.Item(119, 16).Value = 1000
$ws.Cells.Item(119, 17).Value = 1
$ws.Cells.Item(119, 18).Value = "Hortaliza"

# The three rows that were shifted from 121-123 to 122-124 change their
# "Fecha" value from 44250 to 44223.
$ws.Cells.Item(122, 4).Value = 44223
$ws.Cells.Item(123, 4).Value = 44223
$ws.Cells.Item(124, 4).Value = 44223
